# Auto-generated script applying scheduled market-data refresh to Ultima_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 926.5
$ws.Cells.Item(28, 10).Value = 1503.6666
$ws.Cells.Item(28, 12).Value = 1503.6666
$ws.Cells.Item(28, 14).Value = -2473.6666

$ws.Cells.Item(92, 8).Value = 913.1515000000001
$ws.Cells.Item(92, 9).Value = 844.96
$ws.Cells.Item(92, 10).Value = 1126.25
$ws.Cells.Item(92, 11).Value = 844.96
$ws.Cells.Item(92, 12).Value = 1126.25
$ws.Cells.Item(92, 13).Value = 403.04
$ws.Cells.Item(92, 14).Value = -3622.25

$ws.Cells.Item(127, 8).Value = 772.6667
$ws.Cells.Item(127, 9).Value = 469.57144
$ws.Cells.Item(127, 10).Value = 897.4706
$ws.Cells.Item(127, 11).Value = 1408.71432
$ws.Cells.Item(127, 12).Value = 2692.4118
$ws.Cells.Item(127, 13).Value = 3551.28568
$ws.Cells.Item(127, 14).Value = -12612.4118

$ws.Cells.Item(132, 8).Value = 2224.0134
$ws.Cells.Item(132, 9).Value = 1663.4348
$ws.Cells.Item(132, 10).Value = 9960
$ws.Cells.Item(132, 11).Value = 4990.3044
$ws.Cells.Item(132, 12).Value = 29880
$ws.Cells.Item(132, 13).Value = -2460.3044
$ws.Cells.Item(132, 14).Value = -34940

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2605
$ws.Cells.Item(102, 9).Value = 2605
$ws.Cells.Item(102, 11).Value = 2605
$ws.Cells.Item(102, 13).Value = -983

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1524.9584
$ws.Cells.Item(20, 9).Value = 1378.4445
$ws.Cells.Item(20, 10).Value = 1964.5
$ws.Cells.Item(20, 11).Value = 1378.4445
$ws.Cells.Item(20, 12).Value = 1964.5
$ws.Cells.Item(20, 13).Value = -1131.4445
$ws.Cells.Item(20, 14).Value = -2458.5

$ws.Cells.Item(75, 8).Value = 48773.715
$ws.Cells.Item(75, 9).Value = 40000
$ws.Cells.Item(75, 10).Value = 50236
$ws.Cells.Item(75, 11).Value = 40000
$ws.Cells.Item(75, 12).Value = 50236
$ws.Cells.Item(75, 13).Value = -39064
$ws.Cells.Item(75, 14).Value = -52108

$ws.Cells.Item(78, 8).Value = 48773.715
$ws.Cells.Item(78, 9).Value = 40000
$ws.Cells.Item(78, 10).Value = 50236
$ws.Cells.Item(78, 11).Value = 120000
$ws.Cells.Item(78, 12).Value = 150708
$ws.Cells.Item(78, 13).Value = -115320
$ws.Cells.Item(78, 14).Value = -160068

$ws.Cells.Item(99, 8).Value = 1455
$ws.Cells.Item(99, 9).Value = 1455
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 1455
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 43
$ws.Cells.Item(99, 14).ClearContents()

$ws.Cells.Item(105, 8).Value = 4994.737
$ws.Cells.Item(105, 9).Value = 4475
$ws.Cells.Item(105, 10).Value = 5133.3335
$ws.Cells.Item(105, 11).Value = 4475
$ws.Cells.Item(105, 12).Value = 5133.3335
$ws.Cells.Item(105, 13).Value = -2728
$ws.Cells.Item(105, 14).Value = -8627.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 82.75
$ws.Cells.Item(7, 9).Value = 78.71429000000001
$ws.Cells.Item(7, 10).Value = 88.40000000000001
$ws.Cells.Item(7, 11).Value = 78.71429000000001
$ws.Cells.Item(7, 12).Value = 88.40000000000001
$ws.Cells.Item(7, 13).Value = 34.28570999999999
$ws.Cells.Item(7, 14).Value = -314.4

$ws.Cells.Item(31, 8).Value = 4904978.5
$ws.Cells.Item(31, 9).Value = 3414.4082
$ws.Cells.Item(31, 10).Value = 17545854
$ws.Cells.Item(31, 11).Value = 3414.4082
$ws.Cells.Item(31, 12).Value = 17545854
$ws.Cells.Item(31, 13).Value = -3119.4082
$ws.Cells.Item(31, 14).Value = -17546444

$ws.Cells.Item(34, 8).Value = 4904978.5
$ws.Cells.Item(34, 9).Value = 3414.4082
$ws.Cells.Item(34, 10).Value = 17545854
$ws.Cells.Item(34, 11).Value = 3414.4082
$ws.Cells.Item(34, 12).Value = 17545854
$ws.Cells.Item(34, 13).Value = -3212.4082
$ws.Cells.Item(34, 14).Value = -17546258

$ws.Cells.Item(62, 8).Value = 2196
$ws.Cells.Item(62, 9).Value = 2196
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 2196
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -1572
$ws.Cells.Item(62, 14).ClearContents()

$ws.Cells.Item(65, 8).Value = 2196
$ws.Cells.Item(65, 9).Value = 2196
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 10980
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -7860
$ws.Cells.Item(65, 14).ClearContents()

$ws.Cells.Item(70, 8).Value = 41652
$ws.Cells.Item(70, 10).Value = 41652
$ws.Cells.Item(70, 12).Value = 41652
$ws.Cells.Item(70, 14).Value = -42282

$ws.Cells.Item(73, 8).Value = 41652
$ws.Cells.Item(73, 10).Value = 41652
$ws.Cells.Item(73, 12).Value = 41652
$ws.Cells.Item(73, 14).Value = -43836

$ws.Cells.Item(105, 8).Value = 1441.6
$ws.Cells.Item(105, 9).Value = 1417.2307
$ws.Cells.Item(105, 10).Value = 1600
$ws.Cells.Item(105, 11).Value = 1417.2307
$ws.Cells.Item(105, 12).Value = 1600
$ws.Cells.Item(105, 13).Value = 329.7692999999999
$ws.Cells.Item(105, 14).Value = -5094

$ws.Cells.Item(140, 8).Value = 45551.285
$ws.Cells.Item(140, 10).Value = 45551.285
$ws.Cells.Item(140, 12).Value = 45551.285
$ws.Cells.Item(140, 14).Value = -55911.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 603951.9399999999
$ws.Cells.Item(2, 9).Value = 101.333336
$ws.Cells.Item(2, 10).Value = 1262698
$ws.Cells.Item(2, 11).Value = 608.000016
$ws.Cells.Item(2, 12).Value = 7576188
$ws.Cells.Item(2, 13).Value = -495.000016
$ws.Cells.Item(2, 14).Value = -7576414

$ws.Cells.Item(38, 8).Value = 168.46666
$ws.Cells.Item(38, 9).Value = 308.83334
$ws.Cells.Item(38, 10).Value = 74.888885
$ws.Cells.Item(38, 11).Value = 926.5000200000001
$ws.Cells.Item(38, 12).Value = 224.666655
$ws.Cells.Item(38, 13).Value = -579.5000200000001
$ws.Cells.Item(38, 14).Value = -918.666655

$ws.Cells.Item(97, 8).Value = 35002.332
$ws.Cells.Item(97, 9).Value = 51503.5
$ws.Cells.Item(97, 10).Value = 2000
$ws.Cells.Item(97, 11).Value = 154510.5
$ws.Cells.Item(97, 12).Value = 6000
$ws.Cells.Item(97, 13).Value = -154014.5
$ws.Cells.Item(97, 14).Value = -6992

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).ClearContents()

$ws.Cells.Item(70, 8).Value = 12488
$ws.Cells.Item(70, 9).Value = 31385.715
$ws.Cells.Item(70, 10).Value = 5138.8887
$ws.Cells.Item(70, 11).Value = 31385.715
$ws.Cells.Item(70, 12).Value = 5138.8887
$ws.Cells.Item(70, 13).Value = -31115.715
$ws.Cells.Item(70, 14).Value = -5678.8887

$ws.Cells.Item(73, 8).Value = 12488
$ws.Cells.Item(73, 9).Value = 31385.715
$ws.Cells.Item(73, 10).Value = 5138.8887
$ws.Cells.Item(73, 11).Value = 31385.715
$ws.Cells.Item(73, 12).Value = 5138.8887
$ws.Cells.Item(73, 13).Value = -30449.715
$ws.Cells.Item(73, 14).Value = -7010.8887

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 903.4231
$ws.Cells.Item(16, 9).Value = 819.5599999999999
$ws.Cells.Item(16, 10).Value = 3000
$ws.Cells.Item(16, 11).Value = 819.5599999999999
$ws.Cells.Item(16, 12).Value = 3000
$ws.Cells.Item(16, 13).Value = -649.5599999999999
$ws.Cells.Item(16, 14).Value = -3340

$ws.Cells.Item(61, 8).Value = 1498.6923
$ws.Cells.Item(61, 9).Value = 1353
$ws.Cells.Item(61, 10).Value = 2300
$ws.Cells.Item(61, 11).Value = 1353
$ws.Cells.Item(61, 12).Value = 2300
$ws.Cells.Item(61, 13).Value = -1151
$ws.Cells.Item(61, 14).Value = -2704

$ws.Cells.Item(113, 8).Value = 1498.6923
$ws.Cells.Item(113, 9).Value = 1353
$ws.Cells.Item(113, 10).Value = 2300
$ws.Cells.Item(113, 11).Value = 1353
$ws.Cells.Item(113, 12).Value = 2300
$ws.Cells.Item(113, 13).Value = 817
$ws.Cells.Item(113, 14).Value = -6640

$ws.Cells.Item(122, 8).Value = 3691.3157
$ws.Cells.Item(122, 9).Value = 3854.8
$ws.Cells.Item(122, 10).Value = 3376.923
$ws.Cells.Item(122, 11).Value = 11564.4
$ws.Cells.Item(122, 12).Value = 10130.769
$ws.Cells.Item(122, 13).Value = -9114.400000000001
$ws.Cells.Item(122, 14).Value = -15030.769

$ws.Cells.Item(136, 8).Value = 10423111
$ws.Cells.Item(136, 9).Value = 17243666
$ws.Cells.Item(136, 10).Value = 12789.737
$ws.Cells.Item(136, 11).Value = 51730998
$ws.Cells.Item(136, 12).Value = 38369.211
$ws.Cells.Item(136, 13).Value = -51728448
$ws.Cells.Item(136, 14).Value = -43469.211

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2854.3928
$ws.Cells.Item(96, 9).Value = 2215.0667
$ws.Cells.Item(96, 10).Value = 3592.077
$ws.Cells.Item(96, 11).Value = 2215.0667
$ws.Cells.Item(96, 12).Value = 3592.077
$ws.Cells.Item(96, 13).Value = -842.0666999999999
$ws.Cells.Item(96, 14).Value = -6338.077

$ws.Cells.Item(132, 8).Value = 1476.0968
$ws.Cells.Item(132, 9).Value = 1226.7
$ws.Cells.Item(132, 10).Value = 2515.25
$ws.Cells.Item(132, 11).Value = 3680.1
$ws.Cells.Item(132, 12).Value = 7545.75
$ws.Cells.Item(132, 13).Value = -1150.1
$ws.Cells.Item(132, 14).Value = -12605.75
